# "Working on crosshair and smoke power"
# Add a new time-log entry (row 11) for Monday, 2025-05-05, documenting
# work done on crosshair logic.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting (number formats, styles) of the last existing data
# row (row 10) down onto the new row 11 before writing values, so the new
# row's cells pick up the same styles (date / time number formats) as the
# rest of the log instead of receiving default formatting.
$ws.Range("A10:F10").Copy()
$ws.Range("A11:F11").PasteSpecial(-4122)
$null = $ws.Range("A11:F11").ClearContents()

# Day
$ws.Range("A11").Value = "Monday"
# Date
$ws.Range("B11").Value = 45782
# From
$ws.Range("C11").Value = 0.35416666666666669
# Until
$ws.Range("D11").Value = 0.4375
# Time spent
$ws.Range("E11").Formula = "=D11-C11"
# Realised
$ws.Range("F11").Value = "Figuring out crosshair logic"

# Match the selection left behind by the edit (Excel leaves the cursor on
# the last touched cell).
$null = $ws.Range("F11").Select()
